# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp (A1): 13:10 -> 13:40
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 13:40"

# Row 4 - Estados Unidos (A4 unchanged)
$ws.Range("B4").Value = 1793653
$ws.Range("C4").Value = 123
$ws.Range("D4").Value = 519611
$ws.Range("E4").Value = 1169495
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 104547

# Row 23 - Catar (A23 unchanged)
$ws.Range("B23").Value = 55262
$ws.Range("C23").Value = 2355
$ws.Range("D23").Value = 25839
$ws.Range("E23").Value = 29387

# Rows 35 & 36 - Kuwait overtakes Indonesia in the ranking.
# Row 35 becomes Kuwait with freshly updated figures, row 36 becomes
# Indonesia carrying the figures Kuwait's former row 35 slot had before.
$ws.Range("A35").Value = "Kuwait"
$ws.Range("B35").Value = 26192
$ws.Range("C35").Value = 1008
$ws.Range("D35").Value = 10156
$ws.Range("E35").Value = 15831
$ws.Range("G35").Value = 11
$ws.Range("H35").Value = 205

$ws.Range("A36").Value = "Indonesia"
$ws.Range("B36").Value = 25773
$ws.Range("C36").Value = 557
$ws.Range("D36").Value = 7015
$ws.Range("E36").Value = 17185
$ws.Range("G36").Value = 53
$ws.Range("H36").Value = 1573

# Row 66 - Australia (A66 unchanged)
$ws.Range("B66").Value = 7185
$ws.Range("C66").Value = 12
$ws.Range("D66").Value = 6606

# Row 99 - Maldivas (A99 unchanged)
$ws.Range("B99").Value = 1633
$ws.Range("C99").Value = 42
$ws.Range("E99").Value = 1398
